# Updates betexplorer-scraped Turkey Super Lig 2023-2024 rows:
#  - A batch of match rows had been written with home/away (and their odds)
#    mixed up between two adjacent fixtures sharing the same matchday
#    timestamp; this swaps columns F:V (home..url) back between the two
#    rows for every affected pair (A:E - index/country/tournament/season/
#    date - stay put).
#  - Rows 97-99 needed a 3-way rotation instead of a plain swap.
#  - Two newly scraped fixtures are appended as rows 132 and 133.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns F..V (6..22) hold home/away teams, goals, odds, timestamps, url.
$firstCol = 6
$lastCol = 22

function Swap-Rows($rowA, $rowB) {
    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $cellA = $ws.Cells.Item($rowA, $col)
        $cellB = $ws.Cells.Item($rowB, $col)
        $valA = $cellA.Value()
        $valB = $cellB.Value()
        $cellA.Value = $valB
        $cellB.Value = $valA
    }
}

function Rotate-Rows($rowA, $rowB, $rowC) {
    # new A <- old B, new B <- old C, new C <- old A
    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $cellA = $ws.Cells.Item($rowA, $col)
        $cellB = $ws.Cells.Item($rowB, $col)
        $cellC = $ws.Cells.Item($rowC, $col)
        $valA = $cellA.Value()
        $valB = $cellB.Value()
        $valC = $cellC.Value()
        $cellA.Value = $valB
        $cellB.Value = $valC
        $cellC.Value = $valA
    }
}

# Plain row-pair swaps (home/away fixtures that had been transposed).
Swap-Rows 8 9
Swap-Rows 18 19
Swap-Rows 29 30
Swap-Rows 31 32
Swap-Rows 42 43
Swap-Rows 44 45
Swap-Rows 51 52
Swap-Rows 63 64
Swap-Rows 79 80
Swap-Rows 89 90
Swap-Rows 106 107
Swap-Rows 127 128

# Three-way rotation.
Rotate-Rows 97 98 99

# Append the two newly scraped fixtures as rows 132 and 133, matching the
# formatting of the last existing data row (bold/bordered/centered index in
# column A, datetime-formatted match date in column E).
$ws.Range("A131").Copy()
$ws.Range("A132").PasteSpecial(-4122)
$ws.Range("A133").PasteSpecial(-4122)

$ws.Range("E131").Copy()
$ws.Range("E132").PasteSpecial(-4122)
$ws.Range("E133").PasteSpecial(-4122)

function Set-Row($r, $idx, $home, $homeGoals, $away, $awayGoals,
                  $homeOpenOdds, $homeOpenDt, $homeCloseOdds, $homeCloseDt,
                  $drawOpenOdds, $drawOpenDt, $drawCloseOdds, $drawCloseDt,
                  $awayOpenOdds, $awayOpenDt, $awayCloseOdds, $awayCloseDt,
                  $url) {
    $ws.Cells.Item($r, 1).Value = $idx
    $ws.Cells.Item($r, 2).Value = "turkey"
    $ws.Cells.Item($r, 3).Value = "super-lig"
    $ws.Cells.Item($r, 4).Value = "2023-2024"
    $ws.Cells.Item($r, 5).Value = 45261.75
    $ws.Cells.Item($r, 6).Value = $home
    $ws.Cells.Item($r, 7).Value = $homeGoals
    $ws.Cells.Item($r, 8).Value = $away
    $ws.Cells.Item($r, 9).Value = $awayGoals
    $ws.Cells.Item($r, 10).Value = $homeOpenOdds
    $ws.Cells.Item($r, 11).Value = $homeOpenDt
    $ws.Cells.Item($r, 12).Value = $homeCloseOdds
    $ws.Cells.Item($r, 13).Value = $homeCloseDt
    $ws.Cells.Item($r, 14).Value = $drawOpenOdds
    $ws.Cells.Item($r, 15).Value = $drawOpenDt
    $ws.Cells.Item($r, 16).Value = $drawCloseOdds
    $ws.Cells.Item($r, 17).Value = $drawCloseDt
    $ws.Cells.Item($r, 18).Value = $awayOpenOdds
    $ws.Cells.Item($r, 19).Value = $awayOpenDt
    $ws.Cells.Item($r, 20).Value = $awayCloseOdds
    $ws.Cells.Item($r, 21).Value = $awayCloseDt
    $ws.Cells.Item($r, 22).Value = $url
}

Set-Row 132 131 "Karagumruk" 3 "Istanbulspor AS" 0 `
    1.78 "28/11/2023 06:42" 1.72 "01/12/2023 17:59" `
    3.82 "28/11/2023 06:42" 3.86 "01/12/2023 17:58" `
    4.59 "28/11/2023 06:42" 5.26 "01/12/2023 17:59" `
    "https://www.betexplorer.com/football/turkey/super-lig/f-karagumruk-istanbulspor-as/fuPIK0uh/"

Set-Row 133 132 "Hatayspor" 3 "Antalyaspor" 3 `
    2.81 "26/11/2023 14:13" 3.47 "01/12/2023 17:58" `
    3.42 "26/11/2023 14:13" 3.39 "01/12/2023 17:40" `
    2.57 "26/11/2023 14:13" 2.25 "01/12/2023 17:58" `
    "https://www.betexplorer.com/football/turkey/super-lig/hatayspor-antalyaspor/Ucid3aAH/"

Write-Host "done"
